$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Participant Identifier Scheme")

# Refresh the AutoFilter first (while the sheet still has 72 data rows), so the
# subsequent row insertion grows it to the correct A1:J72 footprint (the very
# last data row, 73, stays outside the filtered/autofiltered block, matching
# the pre-existing convention in this sheet of the filter trailing dimension
# by one row).
$ws.AutoFilterMode = $false
$ws.Range("A1:J72").AutoFilter()

# Insert a new row before row 17 (shifts existing rows 17.. down by one).
# The original row 17 (DK:CPR) becomes row 18, etc.
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with the UBLPE / 0193 participant identifier scheme.
$ws.Cells.Item(17, 1).Value = "UBLPE"
$ws.Cells.Item(17, 2).Value = "0193"
$ws.Cells.Item(17, 3).Value = "UBL.BE"
$ws.Cells.Item(17, 4).Value = "3"
$ws.Cells.Item(17, 5).Formula = "=FALSE"
$ws.Cells.Item(17, 7).Value = "Maximum 50 characters`n4 Characters fixed length identifying the type `nMaximum 46 characters for the identifier itself"
$ws.Cells.Item(17, 8).Value = "None"
$ws.Rows.Item(17).RowHeight = 45

# Keep the workbook-level _FilterDatabase defined name in sync with the filter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Participant Identifier Scheme!_FilterDatabase") {
        $n.RefersTo = "='Participant Identifier Scheme'!`$A`$1:`$J`$72"
    }
}

# Mirror the author's final selection (they had just finished typing into D17).
$ws.Range("D17").Select()
